$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 corresponds to group "fresheroffcampusdiscussion"
# D11: current_phase 1 -> 2
$ws.Range("D11").Value = 2

# E11: last_action_date
$ws.Range("E11").Value = "2026-02-21T13:49:36.875132+00:00"

# H11: reactions_count 4 -> 6
$ws.Range("H11").Value = 6

# I11: replies_count 0 -> 2
$ws.Range("I11").Value = 2

# L11: reacted_message_ids - append 68192, 68188
$ws.Range("L11").Value = "[68050, 68073, 68060, 68065, 68192, 68188]"

# M11: replied_message_ids - from [] to [68188, 68187]
$ws.Range("M11").Value = "[68188, 68187]"
